# Weekly refresh of the daily "Alcachofa" price records (rows 3-16).
# Each row's Fecha/Volumen/Precio/Unidad/Origen fields are updated in
# place to reflect the newly consolidated week of data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D3").Value = 44425
$ws.Range("J3").Value = 35

$ws.Range("D4").Value = 44460
$ws.Range("K4").Value = 13000
$ws.Range("L4").Value = 13000
$ws.Range("M4").Value = 13000
$ws.Range("P4").Value = 433

$ws.Range("D5").Value = 44453
$ws.Range("I5").Value = 'Primera'
$ws.Range("J5").Value = 50
$ws.Range("K5").Value = 12000
$ws.Range("L5").Value = 12000
$ws.Range("M5").Value = 12000
$ws.Range("N5").Value = '$/caja 30 unidades'
$ws.Range("P5").Value = 400
$ws.Range("Q5").Value = 30

$ws.Range("D7").Value = 44435
$ws.Range("J7").Value = 25
$ws.Range("K7").Value = 14000
$ws.Range("L7").Value = 14000
$ws.Range("M7").Value = 14000
$ws.Range("P7").Value = 467

$ws.Range("D8").Value = 44435
$ws.Range("J8").Value = 25
$ws.Range("K8").Value = 14000
$ws.Range("L8").Value = 14000
$ws.Range("M8").Value = 14000
$ws.Range("O8").Value = 'Provincia del Elquí'
$ws.Range("P8").Value = 467

$ws.Range("D9").Value = 44446
$ws.Range("J9").Value = 25
$ws.Range("K9").Value = 14000
$ws.Range("L9").Value = 14000
$ws.Range("M9").Value = 14000
$ws.Range("P9").Value = 467

$ws.Range("D10").Value = 44474
$ws.Range("K10").Value = 10000
$ws.Range("L10").Value = 10000
$ws.Range("M10").Value = 10000
$ws.Range("P10").Value = 333

$ws.Range("D11").Value = 44418
$ws.Range("J11").Value = 30
$ws.Range("K11").Value = 15000
$ws.Range("L11").Value = 15000
$ws.Range("M11").Value = 15000
$ws.Range("O11").Value = 'Provincia de Limarí'
$ws.Range("P11").Value = 500

$ws.Range("D12").Value = 44432
$ws.Range("J12").Value = 25
$ws.Range("K12").Value = 14000
$ws.Range("L12").Value = 14000
$ws.Range("M12").Value = 14000
$ws.Range("O12").Value = 'Provincia del Elquí'
$ws.Range("P12").Value = 467

$ws.Range("D13").Value = 44841
$ws.Range("J13").Value = 45
$ws.Range("K13").Value = 12000
$ws.Range("L13").Value = 12000
$ws.Range("M13").Value = 12000
$ws.Range("P13").Value = 400

$ws.Range("D14").Value = 44841
$ws.Range("I14").Value = 'Segunda'
$ws.Range("J14").Value = 45
$ws.Range("K14").Value = 10000
$ws.Range("L14").Value = 10000
$ws.Range("M14").Value = 10000
$ws.Range("N14").Value = '$/caja 40 unidades'
$ws.Range("P14").Value = 250
$ws.Range("Q14").Value = 40

$ws.Range("D15").Value = 44376
$ws.Range("K15").Value = 18000
$ws.Range("L15").Value = 18000
$ws.Range("M15").Value = 18000
$ws.Range("P15").Value = 600

$ws.Range("D16").Value = 44449
$ws.Range("J16").Value = 45
$ws.Range("K16").Value = 12000
$ws.Range("L16").Value = 12000
$ws.Range("M16").Value = 12000
$ws.Range("O16").Value = 'Provincia de Limarí'
$ws.Range("P16").Value = 400

